$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.735.28"
$ws.Range("E2").Value = "  +2.25%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.938.54"
$ws.Range("E3").Value = "  +0.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.20"
$ws.Range("E5").Value = "  +10.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.14"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -0.16%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.727"
$ws.Range("E9").Value = "  -0.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("E10").Value = "  +4.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000336"
$ws.Range("E11").Value = "  -1.14%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.71"
$ws.Range("E12").Value = "  -0.88%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.571.04"
$ws.Range("E13").Value = "  +0.89%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.39"
$ws.Range("E14").Value = "  -3.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.939.54"
$ws.Range("E15").Value = "  +0.63%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.01"
$ws.Range("E16").Value = "  -1.55%  "

# Row 17
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("E18").Value = "  +7.44%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.83"
$ws.Range("E19").Value = "  -0.81%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.592.13"
$ws.Range("E20").Value = "  +1.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "431.85"
$ws.Range("E21").Value = "  +0.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.38"
$ws.Range("E22").Value = "  -2.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.25"
$ws.Range("E23").Value = "  +18.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.51"
$ws.Range("E24").Value = "  -3.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.35"
$ws.Range("E25").Value = "  +0.53%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.68"
$ws.Range("E26").Value = "  +4.99%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.78"
$ws.Range("E27").Value = "  -4.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.68"
$ws.Range("E28").Value = "  -3.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "704.09"
$ws.Range("E29").Value = "  -1.12%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.27"
$ws.Range("E30").Value = "  -4.05%  "

# Row 31
$ws.Range("E31").Value = "  -1.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "70.29"
$ws.Range("E32").Value = "  +16.12%  "

# Row 33
$ws.Range("E33").Value = "  -2.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.464"
$ws.Range("E34").Value = "  +19.23%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.00"
$ws.Range("E35").Value = "  -3.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.26"
$ws.Range("E36").Value = "  -2.55%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0834"
$ws.Range("E37").Value = "  -2.82%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.149"
$ws.Range("E38").Value = "  +2.97%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("E40").Value = "  -0.02%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0482"
$ws.Range("E41").Value = "  -0.61%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.11"
$ws.Range("E42").Value = "  +6.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.78"
$ws.Range("E43").Value = "  -4.06%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.99"
$ws.Range("E44").Value = "  -4.47%  "

# Row 45
$ws.Range("E45").Value = "  +14.22%  "

# Row 46
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.142"
$ws.Range("E46").Value = "  +1.16%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.38"
$ws.Range("E47").Value = "  +0.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0₆0353"
$ws.Range("E48").Value = "  +4.08%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.32"
$ws.Range("E49").Value = "  -2.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.90"
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("E51").Value = "  -1.56%  "
